$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column H
$ws.Range("H1").Value = "Issued By"

# Fix author's name typo
$ws.Range("C2").Value = "George Orwell"

# Book has been returned - update status
$ws.Range("E2").Value = "Available"

# Clear issue/return dates back to zero (book no longer issued)
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Move the active selection to H2
$ws.Range("H2").Select()
